# [PROS-14794] Template update for KPI-linear_product_length_out_of_store
# to exclude sub_category:Out of Scope
#
# Appends a new exclusion row (row 6) to the KPI Exclusions Template,
# mirroring the existing row 2 pattern:
#   KPI = linear_product_length_out_of_store
#   Exclude1 = product_name          Value1 = Empty; Irrelevant; General Empty
#   Exclude2 = sub_category          Value2 = Out of Scope

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -----------------------------------------------------
$ws.Range("A6").Value2 = "linear_product_length_out_of_store"
$ws.Range("B6").Value2 = "product_name"
$ws.Range("C6").Value2 = "Empty; Irrelevant; General Empty"
$ws.Range("D6").Value2 = "sub_category"
$ws.Range("E6").Value2 = "Out of Scope"

# --- Formatting: match the style already used on the analogous row 2/3 -
# B6 (Exclude1 label) is vertical-centered, like B2/B3.
# C6 (Value1) wraps its text, like C2/C3.
$ws.Range("B6").WrapText = $false
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("C6").WrapText = $true
$ws.Range("C6").VerticalAlignment = -4107

# --- Column widths: nudge to accommodate the new, wider content -------
$ws.Columns.Item(1).ColumnWidth = 29.63
$ws.Columns.Item(2).ColumnWidth = 18.56
$ws.Columns.Item(3).ColumnWidth = 45.23
$ws.Columns.Item(4).ColumnWidth = 11.66

# --- Row height for the new row -----------------------------------
$ws.Range("A6").RowHeight = 14.95

# --- Selection follows the newly entered cell, like Excel would leave it
$ws.Range("A6").Select() | Out-Null
